# Auto-generated Excel COM-interop script to apply numeric updates
# to the Ridill_Profits sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3586239.2
$ws.Range("I132").Value = 761950.8
$ws.Range("J132").Value = 22226542
$ws.Range("K132").Value = 2285852.4
$ws.Range("L132").Value = 66679626
$ws.Range("M132").Value = -2283322.4
$ws.Range("N132").Value = -66684686

$ws.Range("H137").Value = 13899021
$ws.Range("I137").Value = 3473080.8
$ws.Range("J137").Value = 34750900
$ws.Range("K137").Value = 10419242.4
$ws.Range("L137").Value = 104252700
$ws.Range("M137").Value = -10416692.4
$ws.Range("N137").Value = -104257800

$ws.Range("H138").Value = 4009.2825
$ws.Range("I138").Value = 3244.0435
$ws.Range("J138").Value = 4293.161
$ws.Range("K138").Value = 9732.130500000001
$ws.Range("L138").Value = 12879.483
$ws.Range("M138").Value = -4592.130500000001
$ws.Range("N138").Value = -23159.483

$ws.Range("H141").Value = 2133
$ws.Range("I141").Value = 616.871
$ws.Range("J141").Value = 7355.222
$ws.Range("K141").Value = 1850.613
$ws.Range("L141").Value = 22065.666
$ws.Range("M141").Value = 3329.387
$ws.Range("N141").Value = -32425.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2919250.2
$ws.Range("I32").Value = 3463504.8
$ws.Range("J32").Value = 49545.453
$ws.Range("K32").Value = 3463504.8
$ws.Range("L32").Value = 49545.453
$ws.Range("M32").Value = -3463217.8
$ws.Range("N32").Value = -50119.453

$ws.Range("H45").Value = 1854.52
$ws.Range("I45").Value = 1134.6818
$ws.Range("K45").Value = 1134.6818
$ws.Range("M45").Value = -757.6818000000001

$ws.Range("H61").Value = 2478073.8
$ws.Range("I61").Value = 1463016
$ws.Range("J61").Value = 6538304.5
$ws.Range("K61").Value = 1463016
$ws.Range("L61").Value = 6538304.5
$ws.Range("M61").Value = -1462804
$ws.Range("N61").Value = -6538728.5

$ws.Range("H74").Value = 13677317
$ws.Range("I74").Value = 1173.96
$ws.Range("J74").Value = 38099000
$ws.Range("K74").Value = 1173.96
$ws.Range("L74").Value = 38099000
$ws.Range("M74").Value = -299.96
$ws.Range("N74").Value = -38100748

$ws.Range("H77").Value = 13677317
$ws.Range("I77").Value = 1173.96
$ws.Range("J77").Value = 38099000
$ws.Range("K77").Value = 5869.8
$ws.Range("L77").Value = 190495000
$ws.Range("M77").Value = -1501.8
$ws.Range("N77").Value = -190503736

$ws.Range("H132").Value = 26672068
$ws.Range("I132").Value = 27618356
$ws.Range("J132").Value = 20836616
$ws.Range("K132").Value = 82855068
$ws.Range("L132").Value = 62509848
$ws.Range("M132").Value = -82852538
$ws.Range("N132").Value = -62514908

$ws.Range("H136").Value = 2478073.8
$ws.Range("I136").Value = 1463016
$ws.Range("J136").Value = 6538304.5
$ws.Range("K136").Value = 4389048
$ws.Range("L136").Value = 19614913.5
$ws.Range("M136").Value = -4386498
$ws.Range("N136").Value = -19620013.5

$ws.Range("H140").Value = 72343
$ws.Range("J140").Value = 72343
$ws.Range("L140").Value = 72343
$ws.Range("N140").Value = -82703

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14428782
$ws.Range("I134").Value = 20756368
$ws.Range("J134").Value = 47905.727
$ws.Range("K134").Value = 62269104
$ws.Range("L134").Value = 143717.181
$ws.Range("M134").Value = -62266569
$ws.Range("N134").Value = -148787.181

$ws.Range("H140").Value = 40650
$ws.Range("J140").Value = 40650
$ws.Range("L140").Value = 40650
$ws.Range("N140").Value = -51010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1571873.4
$ws.Range("I58").Value = 5957.5264
$ws.Range("K58").Value = 5957.5264
$ws.Range("M58").Value = -5754.5264

$ws.Range("H132").Value = 2178.5908
$ws.Range("I132").Value = 1553.8667
$ws.Range("J132").Value = 3517.2856
$ws.Range("K132").Value = 4661.6001
$ws.Range("L132").Value = 10551.8568
$ws.Range("M132").Value = -2131.6001
$ws.Range("N132").Value = -15611.8568

$ws.Range("H134").Value = 1083486.9
$ws.Range("I134").Value = 1613.0741
$ws.Range("K134").Value = 4839.2223
$ws.Range("M134").Value = -2304.2223

$ws.Range("H136").Value = 1571873.4
$ws.Range("I136").Value = 5957.5264
$ws.Range("K136").Value = 17872.5792
$ws.Range("M136").Value = -15322.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1311.6129
$ws.Range("I14").Value = 1311.6129
$ws.Range("K14").Value = 3934.8387
$ws.Range("M14").Value = -3761.8387

$ws.Range("H131").Value = 925.39435
$ws.Range("I131").Value = 312.5
$ws.Range("J131").Value = 1003.2222
$ws.Range("K131").Value = 937.5
$ws.Range("L131").Value = 3009.6666
$ws.Range("M131").Value = 4102.5
$ws.Range("N131").Value = -13089.6666

$ws.Range("H140").Value = 2385.3513
$ws.Range("I140").Value = 2116.9355
$ws.Range("K140").Value = 6350.806500000001
$ws.Range("M140").Value = -1170.806500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4828.5713
$ws.Range("I41").Value = 900
$ws.Range("J41").Value = 6400
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 6400
$ws.Range("M41").Value = -545
$ws.Range("N41").Value = -7110

$ws.Range("H80").Value = 11600.25
$ws.Range("I80").Value = 4987
$ws.Range("J80").Value = 31440
$ws.Range("K80").Value = 4987
$ws.Range("L80").Value = 31440
$ws.Range("M80").Value = -3989
$ws.Range("N80").Value = -33436

$ws.Range("H83").Value = 11600.25
$ws.Range("I83").Value = 4987
$ws.Range("J83").Value = 31440
$ws.Range("K83").Value = 24935
$ws.Range("L83").Value = 157200
$ws.Range("M83").Value = -19943
$ws.Range("N83").Value = -167184

$ws.Range("H99").Value = 6085
$ws.Range("I99").Value = 2613.3333
$ws.Range("J99").Value = 16500
$ws.Range("K99").Value = 2613.3333
$ws.Range("L99").Value = 16500
$ws.Range("M99").Value = -367.3332999999998
$ws.Range("N99").Value = -20992

$ws.Range("H122").Value = 3057.4285
$ws.Range("I122").Value = 2147.0952
$ws.Range("J122").Value = 4422.9287
$ws.Range("K122").Value = 6441.285600000001
$ws.Range("L122").Value = 13268.7861
$ws.Range("M122").Value = -3991.285600000001
$ws.Range("N122").Value = -18168.7861

$ws.Range("H132").Value = 7524665
$ws.Range("I132").Value = 9172512
$ws.Range("J132").Value = 5052894.5
$ws.Range("K132").Value = 27517536
$ws.Range("L132").Value = 15158683.5
$ws.Range("M132").Value = -27515006
$ws.Range("N132").Value = -15163743.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2090.7058
$ws.Range("I40").Value = 2167
$ws.Range("J40").Value = 1950.8334
$ws.Range("K40").Value = 2167
$ws.Range("L40").Value = 1950.8334
$ws.Range("M40").Value = -2031
$ws.Range("N40").Value = -2222.8334

$ws.Range("H81").Value = 27499.5
$ws.Range("I81").Value = 15000
$ws.Range("J81").Value = 39999
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 39999
$ws.Range("M81").Value = -14002
$ws.Range("N81").Value = -41995

$ws.Range("H84").Value = 27499.5
$ws.Range("I84").Value = 15000
$ws.Range("J84").Value = 39999
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 119997
$ws.Range("M84").Value = -40008
$ws.Range("N84").Value = -129981

$ws.Range("H135").Value = 41106.5
$ws.Range("J135").Value = 41106.5
$ws.Range("L135").Value = 41106.5
$ws.Range("N135").Value = -51246.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 723754.9399999999
$ws.Range("I132").Value = 2013.9122
$ws.Range("J132").Value = 5294781.5
$ws.Range("K132").Value = 6041.7366
$ws.Range("L132").Value = 15884344.5
$ws.Range("M132").Value = -3511.7366
$ws.Range("N132").Value = -15889404.5

$ws.Range("H136").Value = 1175.3334
$ws.Range("I136").Value = 1031.1875
$ws.Range("J136").Value = 1834.2858
$ws.Range("K136").Value = 3093.5625
$ws.Range("L136").Value = 5502.857400000001
$ws.Range("M136").Value = -543.5625
$ws.Range("N136").Value = -10602.8574
